$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Globo"
$ws.Range("B3").Value = "Inter TV Rural"
$ws.Range("C3").Value = "Agricultura"
$ws.Range("D3").Value = "2025-03-31T19:34"
$ws.Range("E3").Value = "Positivo"
$ws.Range("F3").Value = "Com Nota"
$ws.Range("G3").Value = "Teste"
